# Weekly price-sheet update: insert one new daily observation as row 28
# (Magnum variety, Peru origin, $/malla 25 kilos), pushing the existing
# rows 28-85 down to 29-86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 28; all rows below shift
# down by one (this also grows the sheet's used range to A1:R86).
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new market observation.
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value = 45070
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100112031
$ws.Cells.Item(28, 7).Value = "Poroto verde"
$ws.Cells.Item(28, 8).Value = "Magnum"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 800
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14500
$ws.Cells.Item(28, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value = 580
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"
